$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has 10 data rows (rows 2-11). The updated data has 15 data
# rows (rows 2-16), so insert 5 new blank rows right after row 11.
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows.Item(12).Insert()
}

$allRows = @(
    ,@("ECs", "Lpl", "Vldlr", "ECs", 3, 1, 36.183528, 108.550584, 0.08181008990825493, 0.08181008990825493, 2, 0.6666666666666666, 0.3081963333333333, 0.924589, 0.09210955608663024, 0.09210955608663024, 11.151630656664, 100.364675909976, 0.00753549106485667, 0.00753549106485667)
    ,@("ECs", "Lpl", "Vldlr", "FAPs", 3, 1, 36.183528, 108.550584, 0.08181008990825493, 0.08181008990825493, 3, 1, 2.207039333333333, 6.621118, 0.65961009678592, 0.6596100967859201, 79.858469514768, 718.726225632912, 0.05396276132244885, 0.05396276132244886)
    ,@("ECs", "Lpl", "Vldlr", "sCs", 3, 1, 36.183528, 108.550584, 0.08181008990825493, 0.08181008990825493, 3, 1, 0.83074, 2.49222, 0.2482803471274497, 0.2482803471274497, 30.05910405072, 270.53193645648, 0.0203118375209494, 0.02031183752094941)
    ,@("FAPs", "Lpl", "Vldlr", "ECs", 3, 1, 62.232638, 186.697914, 0.1407065035230363, 0.1407065035230363, 2, 0.6666666666666666, 0.3081963333333333, 0.924589, 0.09210955608663024, 0.09210955608663024, 19.17987084526067, 172.618837607346, 0.01296041357800875, 0.01296041357800875)
    ,@("FAPs", "Lpl", "Vldlr", "FAPs", 3, 1, 62.232638, 186.697914, 0.1407065035230363, 0.1407065035230363, 3, 1, 2.207039333333333, 6.621118, 0.65961009678592, 0.6596100967859201, 137.3498798830947, 1236.148918947852, 0.09281143040723835, 0.09281143040723837)
    ,@("FAPs", "Lpl", "Vldlr", "sCs", 3, 1, 62.232638, 186.697914, 0.1407065035230363, 0.1407065035230363, 3, 1, 0.83074, 2.49222, 0.2482803471274497, 0.2482803471274497, 51.69914169212, 465.29227522908, 0.03493465953778917, 0.03493465953778918)
    ,@("M1", "Lpl", "Vldlr", "ECs", 3, 1, 203.1252746666667, 609.375824, 0.4592613794630246, 0.4592613794630246, 2, 0.6666666666666666, 0.3081963333333333, 0.924589, 0.09210955608663024, 0.09210955608663024, 62.60246485959289, 563.4221837363359, 0.04230236179007264, 0.04230236179007264)
    ,@("M1", "Lpl", "Vldlr", "FAPs", 3, 1, 203.1252746666667, 609.375824, 0.4592613794630246, 0.4592613794630246, 3, 1, 2.207039333333333, 6.621118, 0.65961009678592, 0.6596100967859201, 448.3054707834702, 4034.749237051232, 0.3029334429576408, 0.3029334429576408)
    ,@("M1", "Lpl", "Vldlr", "sCs", 3, 1, 203.1252746666667, 609.375824, 0.4592613794630246, 0.4592613794630246, 3, 1, 0.83074, 2.49222, 0.2482803471274497, 0.2482803471274497, 168.7442906765867, 1518.69861608928, 0.1140255747153112, 0.1140255747153112)
    ,@("M2", "Lpl", "Vldlr", "ECs", 3, 1, 135.170856, 405.512568, 0.3056180669374137, 0.3056180669374137, 2, 0.6666666666666666, 0.3081963333333333, 0.924589, 0.09210955608663024, 0.09210955608663024, 41.659162192728, 374.932459734552, 0.02815034447765922, 0.02815034447765922)
    ,@("M2", "Lpl", "Vldlr", "FAPs", 3, 1, 135.170856, 405.512568, 0.3056180669374137, 0.3056180669374137, 3, 1, 2.207039333333333, 6.621118, 0.65961009678592, 0.6596100967859201, 298.327395912336, 2684.946563211024, 0.2015887627121132, 0.2015887627121132)
    ,@("M2", "Lpl", "Vldlr", "sCs", 3, 1, 135.170856, 405.512568, 0.3056180669374137, 0.3056180669374137, 3, 1, 0.83074, 2.49222, 0.2482803471274497, 0.2482803471274497, 112.29183691344, 1010.62653222096, 0.07587895974764122, 0.07587895974764124)
    ,@("sCs", "Lpl", "Vldlr", "ECs", 3, 1, 5.574566, 16.723698, 0.01260396016827052, 0.01260396016827052, 2, 0.6666666666666666, 0.3081963333333333, 0.924589, 0.09210955608663024, 0.09210955608663024, 1.718060801124667, 15.462547210122, 0.001160945176032967, 0.001160945176032967)
    ,@("sCs", "Lpl", "Vldlr", "FAPs", 3, 1, 5.574566, 16.723698, 0.01260396016827052, 0.01260396016827052, 3, 1, 2.207039333333333, 6.621118, 0.65961009678592, 0.6596100967859201, 12.30328642826267, 110.729577854364, 0.008313699386478798, 0.008313699386478798)
    ,@("sCs", "Lpl", "Vldlr", "sCs", 3, 1, 5.574566, 16.723698, 0.01260396016827052, 0.01260396016827052, 3, 1, 0.83074, 2.49222, 0.2482803471274497, 0.2482803471274497, 4.63101495884, 41.67913462956, 0.003129315605758754, 0.003129315605758754)
)

$startRow = 2
for ($i = 0; $i -lt $allRows.Length; $i++) {
    $rowVals = $allRows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
